$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra "Goal" columns (C, D, E) -- only Goal1 remains in column B.
$ws.Range("C1:E9").ClearContents()

# Update the values for the remaining Goal1 column.
$ws.Range("B2").Value = 232
$ws.Range("B3").Value = 148

# Update the active selection to B4.
$ws.Range("B4").Select() | Out-Null
